$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Clone the formatting (style) of row 22 into the new rows 23-29 first,
#     so every new cell already carries the correct style index (s="4"/"6"/"7"). ---
$ws.Range("A22:C22").Copy($ws.Range("A23:C23"))
$ws.Range("A22:C22").Copy($ws.Range("A24:C24"))
$ws.Range("A22:C22").Copy($ws.Range("A25:C25"))
$ws.Range("A22:C22").Copy($ws.Range("A26:C26"))
$ws.Range("A22:C22").Copy($ws.Range("A27:C27"))
$ws.Range("A22:C22").Copy($ws.Range("A28:C28"))
$ws.Range("A22:C22").Copy($ws.Range("A29:C29"))

# --- Now write the cell values/text in the same order the strings were
#     originally introduced, so new shared-string entries land in the
#     same sequence as the target workbook. ---

# Update column C for the existing rows 17-20 (new "Hankerank..." variants)
$ws.Range("C17").Value = "Hankerank tasks on python(intro)"
$ws.Range("C18").Value = "Hankerank tasks on python(basic data types)"
$ws.Range("C19").Value = "Hankerank tasks on python(strings)"
$ws.Range("C20").Value = "Hankerank tasks on python(strings)"

# C21 keeps referencing the pre-existing "tasks:count word occurrences..." string
$ws.Range("C21").Value = "tasks:count word occurrences(case sensitive) and count word occurrences(case insensitive"

# New rows 23/24 (python tasks)
$ws.Range("A23").Value = 43332
$ws.Range("B23").Value = "python"
$ws.Range("C23").Value = "Hankerank tasks on python(strings)"

# New rows 25-27 (Eid holiday)
$ws.Range("A25").Value = 43334
$ws.Range("B25").Value = "off"
$ws.Range("C25").Value = "Eid Celebration"

$ws.Range("A26").Value = 43335
$ws.Range("B26").Value = "off"
$ws.Range("C26").Value = "Eid Celebration"

$ws.Range("A27").Value = 43336
$ws.Range("B27").Value = "off"
$ws.Range("C27").Value = "Eid Celebration"

# Existing row 22 updated last to its new description
$ws.Range("C22").Value = "tasks: Extract links from a webpage and Files I/O"

# New row 24 description
$ws.Range("A24").Value = 43333
$ws.Range("B24").Value = "python"
$ws.Range("C24").Value = "Hankerank tasks on python(strings) and Exceptions"

# New rows 28/29 (python tasks again)
$ws.Range("A28").Value = 43337
$ws.Range("B28").Value = "python"
$ws.Range("C28").Value = "Hankerank tasks on python(strings)"

$ws.Range("A29").Value = 43338
$ws.Range("B29").Value = "python"
$ws.Range("C29").Value = "Hankerank tasks on python(strings)"

# --- Update the selected cell to match the final view state ---
$ws.Range("C24").Select()
